$d = $word.ActiveDocument

# ------------------------------------------------------------------
# "To do list" bullet list currently reads (all at outline level 1):
#   Databse linking
#   Homescreen component
#   Sign In / Sign Out functionality
#   User component design
#   Upcomming challenges page
#   Sort comments frontend
#   Bugfix sort doesn't work ... (carries the _GoBack bookmark)
#
# Target reads:
#   Databse linking
#   Plan full database interactivity            (indented, level 2)
#     Replicate what you've done for challenges to comments.  (level 3, _GoBack bookmark)
#   Code full database interactivity             (indented, level 2)
#   Homescreen component
#   Sign In / Sign Out functionality
#   User component design
#   Upcomming challenges page
#   Sort comments frontend
#   Bugfix sort doesn't work ... (bookmark removed)
#
# Notes on the two insertion primitives as implemented here:
#   Range.InsertParagraphAfter()  -> original paragraph keeps its content,
#                                     a new EMPTY paragraph appears after it.
#   Range.InsertParagraphBefore() -> the paragraph object stays bound to
#                                     the now-EMPTY paragraph, and the
#                                     original content moves into a new
#                                     paragraph that follows it.
# ------------------------------------------------------------------

# Locate the anchor paragraph ("Homescreen component") by its current text.
$homescreen = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text.TrimEnd([char]13, [char]7) -eq "Homescreen component") {
        $homescreen = $cand
        break
    }
}

# Step 1: repurpose this paragraph for the new "Plan full database
# interactivity" bullet, one level deeper than before.
$homescreen.Range.Text = "Plan full database interactivity"
$homescreen.Range.ListFormat.ListLevelNumber = 2

# Step 2: insert the new "Replicate..." bullet right after it, one level
# deeper still.
$homescreen.Range.InsertParagraphAfter()
$replicate = $homescreen.Next()
$replicate.Range.Text = "Replicate what you" + [char]8217 + "ve done for challenges to comments."
$replicate.Range.ListFormat.ListLevelNumber = 3

# Step 3: insert "Code full database interactivity" after that, back at
# the same indent level as "Plan full database interactivity".
$replicate.Range.InsertParagraphAfter()
$codeInteractivity = $replicate.Next()
$codeInteractivity.Range.Text = "Code full database interactivity"
$codeInteractivity.Range.ListFormat.ListLevelNumber = 2

# ------------------------------------------------------------------
# Steps 4-7: shift the remaining original bullets' text down by one slot.
# ------------------------------------------------------------------
$signInOut = $codeInteractivity.Next()
$signInOut.Range.Text = "Homescreen component"

$userComponent = $signInOut.Next()
$userComponent.Range.Text = "Sign In / Sign Out functionality"

$upcoming = $userComponent.Next()
$upcoming.Range.Text = "User component design"

$sortComments = $upcoming.Next()
$sortComments.Range.Text = "Upcomming challenges page"

# Step 8: the old "Bugfix ..." paragraph - insert a new bullet before it
# carrying the text that used to be in $sortComments ("Sort comments
# frontend"), and strip the _GoBack bookmark off of it (it now belongs to
# the "Replicate..." bullet above).
# InsertParagraphBefore() leaves $bugfix bound to the new EMPTY paragraph,
# while the original "Bugfix ..." content slides into $bugfix.Next().
$bugfix = $sortComments.Next()
$bugfix.Range.InsertParagraphBefore()
$newSortComments = $bugfix
$newSortComments.Range.Text = "Sort comments frontend"
$newSortComments.Range.ListFormat.ListLevelNumber = 1
$bugfix = $newSortComments.Next()

$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# Re-anchor _GoBack as a zero-length bookmark at the end of the
# "Replicate..." bullet's text (i.e. right before its paragraph mark),
# matching how it originally sat right after the run text in the
# "Bugfix..." paragraph.
$replicateTextEnd = $d.Range($replicate.Range.End - 1, $replicate.Range.End - 1)
$d.Bookmarks.Add("_GoBack", $replicateTextEnd)
